# Update computed market-price / leve-profit figures pulled in by the scheduled Sheets runner.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 255.1
$ws.Range("I4").Value = 179.11111
$ws.Range("K4").Value = 179.11111
$ws.Range("M4").Value = -65.11111
# Row 131
$ws.Range("H131").Value = 1181.5
$ws.Range("I131").Value = 1057.8
$ws.Range("K131").Value = 3173.4
$ws.Range("M131").Value = 1866.6
# Row 135
$ws.Range("H135").Value = 71429160
$ws.Range("I135").Value = 26316176
$ws.Range("K135").Value = 236845584
$ws.Range("M135").Value = -236843049
# Row 137
$ws.Range("H137").Value = 2711.8333
$ws.Range("I137").Value = 1674.6666
$ws.Range("K137").Value = 5023.9998
$ws.Range("M137").Value = -2473.9998
# Row 141
$ws.Range("H141").Value = 2583.5
$ws.Range("I141").Value = 2110.1538
$ws.Range("K141").Value = 6330.4614
$ws.Range("M141").Value = -1150.4614

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 381326.22
$ws.Range("I2").Value = 726027.9399999999
$ws.Range("J2").Value = 3795.762
$ws.Range("K2").Value = 726027.9399999999
$ws.Range("L2").Value = 3795.762
$ws.Range("M2").Value = -725914.9399999999
$ws.Range("N2").Value = -4021.762
# Row 5
$ws.Range("H5").Value = 222.72728
$ws.Range("I5").Value = 241
$ws.Range("K5").Value = 241
$ws.Range("M5").Value = -129
# Row 45
$ws.Range("H45").Value = 1629.9
$ws.Range("I45").Value = 894
$ws.Range("J45").Value = 1945.2858
$ws.Range("K45").Value = 894
$ws.Range("L45").Value = 1945.2858
$ws.Range("M45").Value = -517
$ws.Range("N45").Value = -2699.2858
# Row 61
$ws.Range("H61").Value = 1000000000
$ws.Range("I61").Value = 1000000000
$ws.Range("K61").Value = 1000000000
$ws.Range("M61").Value = -999999788
# Row 102
$ws.Range("H102").Value = 6494204
$ws.Range("I102").Value = 7143524.5
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 7143524.5
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = -7141902.5
$ws.Range("N102").Value = -4244
# Row 116
$ws.Range("H116").Value = 381326.22
$ws.Range("I116").Value = 726027.9399999999
$ws.Range("J116").Value = 3795.762
$ws.Range("K116").Value = 726027.9399999999
$ws.Range("L116").Value = 3795.762
$ws.Range("M116").Value = -723733.9399999999
$ws.Range("N116").Value = -8383.762000000001
# Row 122
$ws.Range("H122").Value = 1452.931
$ws.Range("I122").Value = 1563.7894
$ws.Range("K122").Value = 4691.3682
$ws.Range("M122").Value = -2241.3682
# Row 136
$ws.Range("H136").Value = 1000000000
$ws.Range("I136").Value = 1000000000
$ws.Range("K136").Value = 3000000000
$ws.Range("M136").Value = -2999997450

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 381326.22
$ws.Range("I3").Value = 726027.9399999999
$ws.Range("J3").Value = 3795.762
$ws.Range("K3").Value = 726027.9399999999
$ws.Range("L3").Value = 3795.762
$ws.Range("M3").Value = -725913.9399999999
$ws.Range("N3").Value = -4023.762
# Row 4
$ws.Range("H4").Value = 222.72728
$ws.Range("I4").Value = 241
$ws.Range("K4").Value = 241
$ws.Range("M4").Value = -126
# Row 86
$ws.Range("H86").Value = 1796.2778
$ws.Range("I86").Value = 1662.8667
$ws.Range("J86").Value = 2463.3333
$ws.Range("K86").Value = 1662.8667
$ws.Range("L86").Value = 2463.3333
$ws.Range("M86").Value = -539.8667
$ws.Range("N86").Value = -4709.3333
# Row 89
$ws.Range("H89").Value = 1796.2778
$ws.Range("I89").Value = 1662.8667
$ws.Range("J89").Value = 2463.3333
$ws.Range("K89").Value = 8314.333500000001
$ws.Range("L89").Value = 12316.6665
$ws.Range("M89").Value = -2698.333500000001
$ws.Range("N89").Value = -23548.6665
# Row 94
$ws.Range("H94").Value = 323.4
$ws.Range("I94").Value = 241.75
$ws.Range("K94").Value = 241.75
$ws.Range("M94").Value = 209.25
# Row 105
$ws.Range("H105").Value = 4999
$ws.Range("I105").Value = 4999
$ws.Range("K105").Value = 4999
$ws.Range("M105").Value = -3252
# Row 107
$ws.Range("H107").Value = 93216.91
$ws.Range("I107").Value = 1709.6666
$ws.Range("K107").Value = 1709.6666
$ws.Range("M107").Value = 210.3334

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 435693.56
$ws.Range("I16").Value = 544192.9
$ws.Range("K16").Value = 544192.9
$ws.Range("M16").Value = -543905.9
# Row 31
$ws.Range("H31").Value = 8186.6763
$ws.Range("I31").Value = 1853.5555
$ws.Range("K31").Value = 1853.5555
$ws.Range("M31").Value = -1558.5555
# Row 34
$ws.Range("H34").Value = 8186.6763
$ws.Range("I34").Value = 1853.5555
$ws.Range("K34").Value = 1853.5555
$ws.Range("M34").Value = -1651.5555
# Row 58
$ws.Range("H58").Value = 29419410
$ws.Range("J58").Value = 3675.6
$ws.Range("L58").Value = 3675.6
$ws.Range("N58").Value = -4081.6
# Row 86
$ws.Range("H86").Value = 3791.0833
$ws.Range("I86").Value = 3863
$ws.Range("K86").Value = 3863
$ws.Range("M86").Value = -2740
# Row 89
$ws.Range("H89").Value = 3791.0833
$ws.Range("I89").Value = 3863
$ws.Range("K89").Value = 19315
$ws.Range("M89").Value = -13699
# Row 108
$ws.Range("H108").Value = 74997.5
$ws.Range("J108").Value = 74997.5
$ws.Range("L108").Value = 74997.5
$ws.Range("N108").Value = -82677.5
# Row 113
$ws.Range("H113").Value = 435693.56
$ws.Range("I113").Value = 544192.9
$ws.Range("K113").Value = 544192.9
$ws.Range("M113").Value = -542022.9
# Row 122
$ws.Range("H122").Value = 3353.4133
$ws.Range("I122").Value = 3322.1167
$ws.Range("J122").Value = 3478.6
$ws.Range("K122").Value = 9966.3501
$ws.Range("L122").Value = 10435.8
$ws.Range("M122").Value = -7516.3501
$ws.Range("N122").Value = -15335.8
# Row 132
$ws.Range("H132").Value = 52632428
$ws.Range("I132").Value = 71429470
$ws.Range("J132").Value = 698
$ws.Range("K132").Value = 214288410
$ws.Range("L132").Value = 2094
$ws.Range("M132").Value = -214285880
$ws.Range("N132").Value = -7154
# Row 136
$ws.Range("H136").Value = 29419410
$ws.Range("J136").Value = 3675.6
$ws.Range("L136").Value = 11026.8
$ws.Range("N136").Value = -16126.8

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2352.875
$ws.Range("I102").Value = 2219.6667
$ws.Range("K102").Value = 2219.6667
$ws.Range("M102").Value = -597.6667000000002
# Row 122
$ws.Range("H122").Value = 175228
$ws.Range("I122").Value = 302899.25
$ws.Range("J122").Value = 4999.6665
$ws.Range("K122").Value = 908697.75
$ws.Range("L122").Value = 14998.9995
$ws.Range("M122").Value = -906247.75
$ws.Range("N122").Value = -19898.9995
# Row 132
$ws.Range("H132").Value = 9618949
$ws.Range("I132").Value = 17859334
$ws.Range("K132").Value = 53578002
$ws.Range("M132").Value = -53575472

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7071.4375
$ws.Range("I7").Value = 6720.5386
$ws.Range("K7").Value = 6720.5386
$ws.Range("M7").Value = -6608.5386
# Row 40
$ws.Range("H40").Value = 4814.3125
$ws.Range("I40").Value = 5016.357
$ws.Range("K40").Value = 5016.357
$ws.Range("M40").Value = -4880.357
# Row 93
$ws.Range("H93").Value = 867.6923
$ws.Range("I93").Value = 849.6
$ws.Range("K93").Value = 849.6
$ws.Range("M93").Value = 398.4
# Row 126
$ws.Range("H126").Value = 7071.4375
$ws.Range("I126").Value = 6720.5386
$ws.Range("K126").Value = 20161.6158
$ws.Range("M126").Value = -17691.6158
# Row 132
$ws.Range("H132").Value = 17784252
$ws.Range("I132").Value = 26674828
$ws.Range("J132").Value = 3099.3333
$ws.Range("K132").Value = 80024484
$ws.Range("L132").Value = 9297.999899999999
$ws.Range("M132").Value = -80021954
$ws.Range("N132").Value = -14357.9999
# Row 137
$ws.Range("H137").Value = 124000
$ws.Range("I137").Value = 124000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 124000
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -118900
$ws.Range("N137").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 12346.167
$ws.Range("I41").Value = 14000
$ws.Range("K41").Value = 14000
$ws.Range("M41").Value = -13610
# Row 81
$ws.Range("H81").Value = 5559.615
$ws.Range("I81").Value = 5347.875
$ws.Range("J81").Value = 5898.4
$ws.Range("K81").Value = 10695.75
$ws.Range("L81").Value = 11796.8
$ws.Range("M81").Value = -9634.75
$ws.Range("N81").Value = -13918.8
# Row 84
$ws.Range("H84").Value = 5559.615
$ws.Range("I84").Value = 5347.875
$ws.Range("J84").Value = 5898.4
$ws.Range("K84").Value = 53478.75
$ws.Range("L84").Value = 58984
$ws.Range("M84").Value = -48174.75
$ws.Range("N84").Value = -69592
# Row 122
$ws.Range("H122").Value = 3148.0688
$ws.Range("I122").Value = 3493.348
$ws.Range("K122").Value = 10480.044
$ws.Range("M122").Value = -8030.044
# Row 126
$ws.Range("H126").Value = 1708.7778
$ws.Range("I126").Value = 1562.3334
$ws.Range("K126").Value = 4687.0002
$ws.Range("M126").Value = -2217.0002
